# Updated cryptos list on Thu May 18 04:50:34 UTC 2023 with GitHub Actions
#
# Applies the refreshed Price (D) / Volume(1h) (E) figures to the cryptos
# sheet, plus the RenderToken/Algorand row swap (rows 43-44).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Range, $Text) {
    # Many "Price" values look like plain decimals ("6.742", "0.1718", ...).
    # A normal Value assignment lets Excel auto-convert those to numbers,
    # which would corrupt the original text-formatted price strings. Force
    # text entry by switching the cell to the "@" (Text) number format, then
    # restore the cell's original style so no visible formatting change is
    # left behind.
    $origStyle = $Range.Style
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = $origStyle
}

# row => (new Price D value or $null, new Volume(1h) E value or $null)
$priceVolumeUpdates = [ordered]@{
    2  = @("27.340.32",    "  +0.82%  ")
    3  = @("1.825.79",     "  -0.13%  ")
    4  = @($null,          "  -0.03%  ")
    5  = @("314.93",       "  +0.59%  ")
    6  = @($null,          "  +0.00%  ")
    7  = @("0.4480",       "  -2.17%  ")
    8  = @("0.3778",       "  +0.99%  ")
    9  = @("0.07442",      "  +1.68%  ")
    10 = @("0.8896",       "  +3.01%  ")
    11 = @($null,          "  +0.06%  ")
    12 = @("1.819.53",     "  -0.39%  ")
    13 = @("6.742",        "  +0.27%  ")
    14 = @("5.465",        "  +1.71%  ")
    15 = @("93.50",        "  +0.55%  ")
    16 = @("0.07128",      "  +0.43%  ")
    17 = @($null,          "  -0.04%  ")
    18 = @("0.000008793",  "  -0.80%  ")
    20 = @($null,          "  +0.72%  ")
    21 = @("27.349.39",    "  +0.92%  ")
    22 = @("5.394",        "  +3.63%  ")
    23 = @("10.99",        "  -0.31%  ")
    24 = @("2.053.03",     $null)
    25 = @("1.971",        "  -1.67%  ")
    26 = @("151.57",       "  -0.40%  ")
    27 = @("2.327",        "  +4.30%  ")
    28 = @("18.64",        "  +0.59%  ")
    29 = @("5.388",        "  +1.84%  ")
    30 = @("117.90",       "  +0.20%  ")
    31 = @("0.08894",      "  -0.20%  ")
    32 = @("0.7968",       "  +4.09%  ")
    33 = @("1.202",        "  +0.06%  ")
    34 = @("4.608",        "  +2.71%  ")
    35 = @("2.921",        "  -1.60%  ")
    36 = @($null,          "  -0.03%  ")
    37 = @("1.112",        "  +0.65%  ")
    38 = @($null,          "  +0.55%  ")
    39 = @("0.05315",      "  +0.29%  ")
    40 = @("7.299",        "  +1.45%  ")
    41 = @("0.5353",       "  -0.52%  ")
    42 = @($null,          "  -0.27%  ")
    45 = @("8.675",        "  +0.41%  ")
    46 = @("0.5092",       "  -2.53%  ")
    47 = @("10.65",        "  -0.86%  ")
    48 = @($null,          "  +0.52%  ")
    49 = @("105.32",       "  -0.89%  ")
    50 = @($null,          "  -0.01%  ")
    51 = @("0.06412",      "  -0.17%  ")
}

foreach ($row in $priceVolumeUpdates.Keys) {
    $pair = $priceVolumeUpdates[$row]
    $newD = $pair[0]
    $newE = $pair[1]
    if ($null -ne $newD) {
        Set-TextValue $ws.Range("D$row") $newD
    }
    if ($null -ne $newE) {
        $ws.Range("E$row").Value = $newE
    }
}

# Rows 43/44 swap places: RenderToken <-> Algorand.
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D43") "0.1718"
$ws.Range("E43").Value = "  -0.21%  "

$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D44") "2.328"
$ws.Range("E44").Value = "  +16.49%  "
